$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update input values (Paquete 1 block) ---
$ws.Range("E8").Value = 2000
$ws.Range("J9").Value = 0.15
$ws.Range("S8").Value = 10000

# --- New "Importadora" mini-table to the right (rows 8, 10, 12) ---
$ws.Range("X8").Value = "Items"
$ws.Range("Y8").Value = 15
$ws.Range("Z8").NumberFormat = "0%"
$ws.Range("Z8").Value = 0.01

$ws.Range("X10").Value = "Polizas Aparte"
$ws.Range("Z10").Value = 70

$ws.Range("X12").Value = "Con factura"
$ws.Range("Z12").NumberFormat = "0%"
$ws.Range("Z12").Value = 0.13

# --- New labels near U16:U18 ---
$ws.Range("U16").Value = "Importadora"
$ws.Range("U17").Value = "Si nosotros importamos = 4%"
$ws.Range("U18").Value = "Si ellos importan = 15%"

# --- Update input value in the Paquete 2 block ---
$ws.Range("L25").Value = 10000

# --- Update selection to match the recorded view ---
$ws.Range("H8:L8").Select() | Out-Null
